# "Updated Code With Item System" - add the "Is Inserted?" (column D) item-system
# flags for the special attacks that were missing them, and fix the stray
# "yA15" typo left in the Ogre's "Is Written?" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Good Guys table: Paladin's first special attack (ShieldBash) now has its item inserted.
$ws.Range("D7").Value = "y"

# Bad Guys table: every special attack now has its item inserted.
$ws.Range("D13").Value = "y"
$ws.Range("D14").Value = "y"
$ws.Range("D15").Value = "y"
$ws.Range("D16").Value = "y"
$ws.Range("D17").Value = "y"
$ws.Range("D18").Value = "y"
$ws.Range("D19").Value = "y"
$ws.Range("D20").Value = "y"

# Fix the stray "yA15" typo in Ogre's "Is Written?" cell -> should just be "y".
$ws.Range("C17").Value = "y"

# Restore the view: scrolled down so row 13 is at the top, with D20 selected.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
[void]$ws.Range("D20").Select()
